# "Buthan added to the graphs" - swap the E1/F1 header labels ("wind"/"solar")
# on Sheet1 of the palette workbook, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the column headers in E1 and F1.
$e1Value = $ws.Range("E1").Value()
$f1Value = $ws.Range("F1").Value()
$ws.Range("E1").Value = $f1Value
$ws.Range("F1").Value = $e1Value

# Update the selected/active cell to match the saved view state.
$ws.Range("E2").Select()
